$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Supplier<T>" row (row 11, previously blank) first so the new
# shared strings land in the same order as the authored edit.
$ws.Range("E11").Value = "get( )"
$ws.Range("C11").Value = "Supplier< T >"
$ws.Range("D11").Value = "T "

# Normalize the generic-style spacing on the existing Functional Interface
# names (Consumer/Predicate/Runnable/Callable) to match "< T >" style.
$ws.Range("C4").Value = "Consumer< T >"
$ws.Range("C6").Value = "Predicate< T >"
$ws.Range("C7").Value = "Runnable< >"
$ws.Range("C8").Value = "Callable< T >"

# Move the active selection to C10, matching the saved view state.
[void]$ws.Range("C10").Select()
